$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Date / text updates ---
$ws.Range("B19").Value = "20 TK"
$ws.Range("F34").Value = "26.12.2024 payment "
$ws.Range("B1").Value = "25.12.2024"

# --- Numeric updates ---
$ws.Range("C9").Value = 451122
$ws.Range("C10").ClearContents()
$ws.Range("C14").Value = 25
$ws.Range("C16").Value = 40
$ws.Range("C17").Value = 75
$ws.Range("C19").Value = 500
$ws.Range("D19").Value = 19.3
$ws.Range("E22").Value = 166863
$ws.Range("E23").Value = 15869
$ws.Range("E27").ClearContents()
$ws.Range("E34").Value = 200000

# --- View state ---
$ws.Range("E23").Select()

Write-Output "done"
